$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.00000001450340058271589
$ws.Range("G3").Value = 36.98554495056295
$ws.Range("H3").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 11.17004962369128
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 13.55117343019618
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 10.29038269963711
$ws.Range("H6").Value = 23.46430640971666
$ws.Range("C7").Value = 36.98554495056295
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 10.29038269963711
$ws.Range("H7").Value = 8.808503108192442
$ws.Range("I7").Value = 11.17004962369128
$ws.Range("B8").Value = 0.00000001450340058271589
$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 13.55117343019618
$ws.Range("F8").Value = 23.46430640971666
$ws.Range("G8").Value = 8.808503108192442
$ws.Range("I8").Value = 0
$ws.Range("D9").Value = 11.17004962369128
$ws.Range("G9").Value = 11.17004962369128
$ws.Range("H9").Value = 0
